$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row with the new value
$ws.Range("A11").Value = "Add new table Level of permission"

# Apply yellow fill formatting to rows 9, 10, 11 (whole rows) to match customFormat rows
$ws.Range("9:11").Interior.Color = 65535

# Select the new row (A11:XFD11) as the active selection
$ws.Rows(11).Select()
